$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-home the "Стандарты МОТ" answer from row 26 up onto row 25 (next to
#    its "7. Сопоставимость..." heading), matching the original author's
#    formatting, then collapse the now-empty row 26 and the stray trailing
#    blank row 28 out of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("B26").Copy()
$ws.Range("B25").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B25").Value = $ws.Range("B26").Value2

$ws.Range("A26:B26").ClearContents()
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# ---------------------------------------------------------------------------
# 2. Give every section heading a shaded filler cell in column B (to match
#    the shaded heading cell in column A), for every heading row that didn't
#    already carry real content in B.
# ---------------------------------------------------------------------------
$fillerRows = @(1, 5, 11, 15, 18, 22)
foreach ($r in $fillerRows) {
    $ws.Range("A1").Copy()
    $ws.Range("B" + $r).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("B" + $r).Font.Bold = $false
    $ws.Range("B" + $r).Borders.LineStyle = -4142   # xlLineStyleNone
    $ws.Range("B" + $r).WrapText = $false
}

# ---------------------------------------------------------------------------
# 3. Shade the label cells in column A (rows with a field caption) with the
#    same accent fill used on the section headings, keeping their existing
#    border/wrap/font otherwise intact.
# ---------------------------------------------------------------------------
$labelRows = @(2, 3, 4, 6, 7, 8, 9, 10, 12, 13, 14, 16, 17, 19, 20, 21, 23, 24)
foreach ($r in $labelRows) {
    $ws.Range("A1").Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("A" + $r).Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 4. Restore the current selection to the cell the author left active.
# ---------------------------------------------------------------------------
$ws.Range("D3").Select()
